$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 489.35822
$ws.Range("I15").Value = 489.35822
$ws.Range("K15").Value = 1468.07466
$ws.Range("M15").Value = -1299.07466
$ws.Range("H32").Value = 6216.6
$ws.Range("I32").Value = 6244.5
$ws.Range("J32").Value = 6209.625
$ws.Range("K32").Value = 6244.5
$ws.Range("L32").Value = 6209.625
$ws.Range("M32").Value = -5918.5
$ws.Range("N32").Value = -6861.625
$ws.Range("H55").Value = 174.88889
$ws.Range("I55").Value = 130.16667
$ws.Range("J55").Value = 264.33334
$ws.Range("K55").Value = 130.16667
$ws.Range("L55").Value = 264.33334
$ws.Range("M55").Value = 83.83332999999999
$ws.Range("N55").Value = -692.33334
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
$ws.Range("H107").Value = 720.5714
$ws.Range("I107").Value = 581.94116
$ws.Range("J107").Value = 1309.75
$ws.Range("K107").Value = 581.94116
$ws.Range("L107").Value = 1309.75
$ws.Range("M107").Value = 1338.05884
$ws.Range("N107").Value = -5149.75
$ws.Range("H125").Value = 3336.8
$ws.Range("I125").Value = 2600
$ws.Range("J125").Value = 3828
$ws.Range("K125").Value = 23400
$ws.Range("L125").Value = 34452
$ws.Range("M125").Value = -20940
$ws.Range("N125").Value = -39372
$ws.Range("H132").Value = 2210.9795
$ws.Range("I132").Value = 1920.2142
$ws.Range("J132").Value = 3955.5715
$ws.Range("K132").Value = 5760.642599999999
$ws.Range("L132").Value = 11866.7145
$ws.Range("M132").Value = -3230.642599999999
$ws.Range("N132").Value = -16926.7145
$ws.Range("H138").Value = 9805711
$ws.Range("I138").Value = 1305.4667
$ws.Range("J138").Value = 17546032
$ws.Range("K138").Value = 3916.4001
$ws.Range("L138").Value = 52638096
$ws.Range("M138").Value = 1223.5999
$ws.Range("N138").Value = -52648376

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1535.3334
$ws.Range("I2").Value = 1324.75
$ws.Range("K2").Value = 1324.75
$ws.Range("M2").Value = -1211.75
$ws.Range("H32").Value = 18525250
$ws.Range("I32").Value = 22731024
$ws.Range("K32").Value = 22731024
$ws.Range("M32").Value = -22730737
$ws.Range("H45").Value = 2008.3334
$ws.Range("I45").Value = 1721
$ws.Range("K45").Value = 1721
$ws.Range("M45").Value = -1344
$ws.Range("H61").Value = 32261524
$ws.Range("I61").Value = 41669150
$ws.Range("J61").Value = 6799.7144
$ws.Range("K61").Value = 41669150
$ws.Range("L61").Value = 6799.7144
$ws.Range("M61").Value = -41668938
$ws.Range("N61").Value = -7223.7144
$ws.Range("H110").Value = 14976.786
$ws.Range("I110").Value = 17342.39
$ws.Range("K110").Value = 17342.39
$ws.Range("M110").Value = -15297.39
$ws.Range("H116").Value = 1535.3334
$ws.Range("I116").Value = 1324.75
$ws.Range("K116").Value = 1324.75
$ws.Range("M116").Value = 969.25
$ws.Range("H122").Value = 4296.579
$ws.Range("I122").Value = 3650
$ws.Range("J122").Value = 4595
$ws.Range("K122").Value = 10950
$ws.Range("L122").Value = 13785
$ws.Range("M122").Value = -8500
$ws.Range("N122").Value = -18685
$ws.Range("H132").Value = 21741844
$ws.Range("I132").Value = 2769.4146
$ws.Range("K132").Value = 8308.2438
$ws.Range("M132").Value = -5778.2438
$ws.Range("H136").Value = 32261524
$ws.Range("I136").Value = 41669150
$ws.Range("J136").Value = 6799.7144
$ws.Range("K136").Value = 125007450
$ws.Range("L136").Value = 20399.1432
$ws.Range("M136").Value = -125004900
$ws.Range("N136").Value = -25499.1432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1535.3334
$ws.Range("I3").Value = 1324.75
$ws.Range("K3").Value = 1324.75
$ws.Range("M3").Value = -1210.75
$ws.Range("H106").Value = 10057
$ws.Range("J106").Value = 10057
$ws.Range("L106").Value = 10057
$ws.Range("N106").Value = -12581

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 69999
$ws.Range("J70").Value = 69999
$ws.Range("L70").Value = 69999
$ws.Range("N70").Value = -70629
$ws.Range("H73").Value = 69999
$ws.Range("J73").Value = 69999
$ws.Range("L73").Value = 69999
$ws.Range("N73").Value = -72183
$ws.Range("H99").Value = 10660.075
$ws.Range("I99").Value = 11481.417
$ws.Range("J99").Value = 10308.071
$ws.Range("K99").Value = 11481.417
$ws.Range("L99").Value = 10308.071
$ws.Range("M99").Value = -9983.416999999999
$ws.Range("N99").Value = -13304.071
$ws.Range("H122").Value = 1825.8636
$ws.Range("I122").Value = 1763.9375
$ws.Range("J122").Value = 1991
$ws.Range("K122").Value = 5291.8125
$ws.Range("L122").Value = 5973
$ws.Range("M122").Value = -2841.8125
$ws.Range("N122").Value = -10873
$ws.Range("H126").Value = 10660.075
$ws.Range("I126").Value = 11481.417
$ws.Range("J126").Value = 10308.071
$ws.Range("K126").Value = 34444.251
$ws.Range("L126").Value = 30924.213
$ws.Range("M126").Value = -31974.251
$ws.Range("N126").Value = -35864.213
$ws.Range("H132").Value = 3460.7778
$ws.Range("I132").Value = 2762.6316
$ws.Range("K132").Value = 8287.8948
$ws.Range("M132").Value = -5757.8948
$ws.Range("H134").Value = 1462.8667
$ws.Range("I134").Value = 1226.6154
$ws.Range("K134").Value = 3679.8462
$ws.Range("M134").Value = -1144.8462

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 428.2857
$ws.Range("I92").Value = 449.75
$ws.Range("K92").Value = 1349.25
$ws.Range("M92").Value = -101.25
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("H128").Value = 121490.336
$ws.Range("I128").Value = 121490.336
$ws.Range("K128").Value = 364471.008
$ws.Range("M128").Value = -359491.008

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = 0
$ws.Range("H97").Value = 444.2
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 4200
$ws.Range("J102").Value = 4666.6665
$ws.Range("L102").Value = 4666.6665
$ws.Range("N102").Value = -7910.6665
$ws.Range("H109").Value = 33864.832
$ws.Range("J109").Value = 33672.25
$ws.Range("L109").Value = 33672.25
$ws.Range("N109").Value = -35752.25
$ws.Range("H113").Value = 3564.95
$ws.Range("I113").Value = 2574.75
$ws.Range("K113").Value = 2574.75
$ws.Range("M113").Value = -404.75
$ws.Range("H114").Value = 49722
$ws.Range("J114").Value = 49722
$ws.Range("L114").Value = 49722
$ws.Range("N114").Value = -58400
$ws.Range("H121").Value = 77749.5
$ws.Range("J121").Value = 77749.5
$ws.Range("L121").Value = 77749.5
$ws.Range("N121").Value = -81243.5
$ws.Range("H122").Value = 6686.7896
$ws.Range("I122").Value = 3038.2856
$ws.Range("K122").Value = 9114.856800000001
$ws.Range("M122").Value = -6664.856800000001
$ws.Range("H132").Value = 2924.8823
$ws.Range("I132").Value = 2632.652
$ws.Range("K132").Value = 7897.956
$ws.Range("M132").Value = -5367.956

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4103.96
$ws.Range("I7").Value = 4280
$ws.Range("J7").Value = 3839.9
$ws.Range("K7").Value = 4280
$ws.Range("L7").Value = 3839.9
$ws.Range("M7").Value = -4168
$ws.Range("N7").Value = -4063.9
$ws.Range("H40").Value = 4338.222
$ws.Range("I40").Value = 4318
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 4318
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -4182
$ws.Range("N40").Value = -4772
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248
$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240
$ws.Range("H126").Value = 4103.96
$ws.Range("I126").Value = 4280
$ws.Range("J126").Value = 3839.9
$ws.Range("K126").Value = 12840
$ws.Range("L126").Value = 11519.7
$ws.Range("M126").Value = -10370
$ws.Range("N126").Value = -16459.7
$ws.Range("H132").Value = 80002830
$ws.Range("I132").Value = 2790.6667
$ws.Range("K132").Value = 8372.000100000001
$ws.Range("M132").Value = -5842.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 342.2143
$ws.Range("I107").Value = 236.45454
$ws.Range("J107").Value = 730
$ws.Range("K107").Value = 709.3636200000001
$ws.Range("L107").Value = 2190
$ws.Range("M107").Value = 1210.63638
$ws.Range("N107").Value = -6030
$ws.Range("H126").Value = 6642.4287
$ws.Range("I126").Value = 7076.4614
$ws.Range("K126").Value = 21229.3842
$ws.Range("M126").Value = -18759.3842
$ws.Range("H132").Value = 3560.5405
$ws.Range("I132").Value = 3560.5405
$ws.Range("K132").Value = 10681.6215
$ws.Range("M132").Value = -8151.621500000001
